# cfast tg: update radiation figures
#
# 1) Remove the second slide (V1/V2/V3 radiation-view diagram).
# 2) Refresh the cached "datetimeFigureOut" footer field (8/3/2015 -> 8/7/2015)
#    on the slide master and every slide layout.
# 3) Nudge the two horizontal divider connectors on slide 1 to their new
#    vertical position.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Delete the old second slide (rId3 / "V1 V2 V3" diagram slide).
# ---------------------------------------------------------------------------
if ($p.Slides.Count -ge 2) {
    $p.Slides.Item(2).Delete()
}

# ---------------------------------------------------------------------------
# 2) Update the cached date field text wherever it still reads 8/3/2015.
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "8/3/2015") {
                $sh.TextFrame.TextRange.Text = "8/7/2015"
            }
        }
    }
}

Update-DateShapes $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 3) Reposition the two horizontal connectors on slide 1.
#    (EMU target / 12700 landed on an f32 storage boundary, so the literal
#    point values below are chosen to round-trip to the exact target EMU.)
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Name -eq "Straight Connector 4") {
        $sh.Top = 133.4268798828125
    } elseif ($sh.Name -eq "Straight Connector 6") {
        $sh.Top = 333.5924987792969
    }
}
